$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = -0.0772
$ws.Range("E2").ClearContents()
$ws.Range("F2").ClearContents()
$ws.Range("G2").Value = -0.05063478004133451
$ws.Range("H2").Value = -0.05063478004133451
$ws.Range("I2").Value = -0.1571739759132892
$ws.Range("J2").Value = -0.1571739759132892
$ws.Range("K2").Value = -107.65
$ws.Range("L2").Value = -0.3178328904635371
$ws.Range("U2").Value = 105.28
$ws.Range("V2").Value = 0.1125267208208636
$ws.Range("W2").Value = -0.2554138173822992
$ws.Range("X2").Value = 0.06318178548876802
$ws.Range("Y2").Value = -0.3185956028710673
$ws.Range("Z2").Value = 0.589110966598915
$ws.Range("AA2").Value = -0.07621042742932216
$ws.Range("AB2").Value = 0.05565803967333459
$ws.Range("AC2").Value = -0.1339943922123562
$ws.Range("AD2").Value = 193.23
$ws.Range("AE2").Value = 72.37412820915523
$ws.Range("AF2").Value = 265.6041282091552
$ws.Range("AG2").Value = 160.3241282091552
$ws.Range("AH2").Value = 0.2211148979359051
$ws.Range("AI2").Value = 0.4815122197715758
$ws.Range("AJ2").Value = 0.1462912660488096
$ws.Range("AK2").Value = 0.359210085397903
$ws.Range("AL2").Value = 7.948
$ws.Range("AM2").Value = 7.136
$ws.Range("AN2").Value = -13.91445236552171
$ws.Range("AO2").Value = -7.453447408152995
$ws.Range("AP2").Value = -11.54490733845721
$ws.Range("AQ2").Value = -8.301569506726457

# Row 3
$ws.Range("B3").Value = "BEC World Public Company Limited (SET:BEC)"
$ws.Range("D3").Value = -0.177
$ws.Range("E3").ClearContents()
$ws.Range("G3").Value = -0.01923857868020305
$ws.Range("H3").Value = -0.01923857868020305
$ws.Range("I3").Value = -0.08088788756743431
$ws.Range("J3").Value = -0.08088788756743431
$ws.Range("K3").Value = -23.5
$ws.Range("L3").Value = -0.1192893401015228
$ws.Range("O3").Value = 0
$ws.Range("R3").Value = 0
$ws.Range("U3").Value = 82.40000000000001
$ws.Range("V3").Value = 0.1536453477531233
$ws.Range("W3").Value = -0.1220145379023884
$ws.Range("X3").Value = 0.06007189286788144
$ws.Range("Y3").Value = -0.1820864307702698
$ws.Range("Z3").Value = 0.74897946058346
$ws.Range("AA3").Value = -0.06058336639799251
$ws.Range("AB3").Value = 0.0548728760389099
$ws.Range("AC3").Value = -0.1154562424369024
$ws.Range("AD3").Value = 110.9
$ws.Range("AE3").Value = 10.62456925392278
$ws.Range("AF3").Value = 121.5245692539228
$ws.Range("AG3").Value = 39.12456925392279
$ws.Range("AH3").Value = 0.1847370483467209
$ws.Range("AI3").Value = 0.4286209465857125
$ws.Range("AJ3").Value = 0.06799252472769882
$ws.Range("AK3").Value = 0.1945290393861698
$ws.Range("AL3").Value = 3.72
$ws.Range("AM3").Value = 3.271
$ws.Range("AN3").Value = -12.04125950054289
$ws.Range("AO3").Value = -5.134408602150538
$ws.Range("AP3").Value = -4.248053122032875
$ws.Range("AQ3").Value = -5.839192907367777

# Row 4
$ws.Range("B4").Value = "Nation Broadcasting Corporation Public Company Limited (SET:NBC)"
$ws.Range("D4").Value = -0.00321
$ws.Range("F4").ClearContents()
$ws.Range("G4").Value = -0.009811320754716982
$ws.Range("H4").Value = -0.009811320754716982
$ws.Range("I4").Value = -0.0876633140696116
$ws.Range("J4").Value = -0.0876633140696116
$ws.Range("K4").Value = -1.05
$ws.Range("L4").Value = -0.03962264150943397
$ws.Range("U4").Value = 4.37
$ws.Range("V4").Value = 0.1942222222222222
$ws.Range("W4").Value = -0.06862745098039216
$ws.Range("X4").Value = 0.09240379141977971
$ws.Range("Y4").Value = -0.1610312424001719
$ws.Range("Z4").Value = 0.7419770764711187
$ws.Range("AA4").Value = -0.06504416948713991
$ws.Range("AB4").Value = 0.0598374098917507
$ws.Range("AC4").Value = -0.1248815793788906
$ws.Range("AD4").Value = 5.84
$ws.Range("AE4").Value = 20.41538911422354
$ws.Range("AF4").Value = 26.25538911422354
$ws.Range("AG4").Value = 21.88538911422354
$ws.Range("AH4").Value = 0.5385125540217253
$ws.Range("AI4").Value = 0.5946134694069773
$ws.Range("AJ4").Value = 0.4930764278736546
$ws.Range("AK4").Value = 0.550086089428226
$ws.Range("AL4").Value = 0.444
$ws.Range("AM4").Value = 0.434
$ws.Range("AN4").Value = 1.455270371293297
$ws.Range("AO4").Value = -3.693693693693693
$ws.Range("AP4").Value = 5.453623003793555
$ws.Range("AQ4").Value = -3.778801843317972

# Row 5
$ws.Range("B5").Value = "MCOT Public Company Limited (SET:MCOT)"
$ws.Range("D5").Value = -0.124
$ws.Range("G5").Value = -0.04858934169278997
$ws.Range("H5").Value = -0.04858934169278997
$ws.Range("I5").Value = -0.1819190484015333
$ws.Range("J5").Value = -0.1819190484015333
$ws.Range("K5").Value = -54.6
$ws.Range("L5").Value = -0.8557993730407524
$ws.Range("U5").Value = 16.1
$ws.Range("V5").Value = 0.165979381443299
$ws.Range("W5").Value = -0.4428223844282239
$ws.Range("X5").Value = 0.06348982638630897
$ws.Range("Y5").Value = -0.5063122108145328
$ws.Range("Z5").Value = 0.4803053123862313
$ws.Range("AA5").Value = -0.08737668537150441
$ws.Range("AB5").Value = 0.05573051967431746
$ws.Range("AC5").Value = -0.1431072050458219
$ws.Range("AD5").Value = 9.390000000000001
$ws.Range("AE5").Value = 22.23217644008913
$ws.Range("AF5").Value = 31.62217644008913
$ws.Range("AG5").Value = 15.52217644008913
$ws.Range("AH5").Value = 0.2458532215462736
$ws.Range("AI5").Value = 0.332088360319981
$ws.Range("AJ5").Value = 0.1379477088976651
$ws.Range("AK5").Value = 0.1961798466431523
$ws.Range("AL5").Value = 0.074
$ws.Range("AM5").Value = -0.279
$ws.Range("AN5").Value = 5.335227272727273
$ws.Range("AO5").Value = -175.6756756756757
$ws.Range("AP5").Value = 8.819418431868822
$ws.Range("AQ5").Value = 46.59498207885305

# Row 6
$ws.Range("B6").Value = "Mono Next Public Company Limited (SET:MONO)"
$ws.Range("D6").Value = -0.0304
$ws.Range("E6").ClearContents()
$ws.Range("G6").Value = -0.1945525291828794
$ws.Range("H6").Value = -0.1945525291828794
$ws.Range("I6").Value = -0.4546770171242015
$ws.Range("J6").Value = -0.4546770171242015
$ws.Range("K6").Value = -28.5
$ws.Range("L6").Value = -0.5544747081712063
$ws.Range("O6").Value = 0
$ws.Range("R6").Value = 0
$ws.Range("U6").Value = 2.41
$ws.Range("V6").Value = 0.008613295210864904
$ws.Range("W6").Value = -0.3888130968622101
$ws.Range("X6").Value = 0.06287374459122705
$ws.Range("Y6").Value = -0.4516868414534372
$ws.Range("Z6").Value = 0.3585329610774666
$ws.Range("AA6").Value = -0.16301669728341
$ws.Range("AB6").Value = 0.05558555967235172
$ws.Range("AC6").Value = -0.2186022569557617
$ws.Range("AD6").Value = 67.09999999999999
$ws.Range("AE6").Value = 19.10199340091979
$ws.Range("AF6").Value = 86.20199340091978
$ws.Range("AG6").Value = 83.79199340091978
$ws.Range("AH6").Value = 0.2355232893676998
$ws.Range("AI6").Value = 0.6697797844699407
$ws.Range("AJ6").Value = 0.2304561016791296
$ws.Range("AK6").Value = 0.6634782708268625
$ws.Range("AL6").Value = 3.71
$ws.Range("AM6").Value = 3.71
$ws.Range("AN6").Value = -6.421052631578947
$ws.Range("AO6").Value = -6.873315363881401
$ws.Range("AP6").Value = -8.018372574250698
$ws.Range("AQ6").Value = -6.873315363881401
